# semana 41 de 2025
# Adds a new weekly column (AR) with the "41" header label and the
# week-41 counts for every UPGD row, mirroring the existing columns
# D:AQ (weeks 1-40). Also backfills AQ44, which was missing its week-40
# value (0) in the source row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: week number 41 in row 1 (stored as text, like the other week
# headers in D1:AQ1, via the leading single-quote so Excel keeps "41"
# as text instead of coercing it to a number).
$ws.Range("AR1").Value = "'41"

# Week-41 case counts per UPGD row.
$ws.Range("AR2").Value = 0
$ws.Range("AR5").Value = 0
$ws.Range("AR6").Value = 2
$ws.Range("AR7").Value = 0
$ws.Range("AR8").Value = 0
$ws.Range("AR9").Value = 0
$ws.Range("AR11").Value = 0
$ws.Range("AR13").Value = 0
$ws.Range("AR14").Value = 0
$ws.Range("AR16").Value = 0
$ws.Range("AR17").Value = 0
$ws.Range("AR22").Value = 0
$ws.Range("AR23").Value = 0
$ws.Range("AR24").Value = 0
$ws.Range("AR25").Value = 0
$ws.Range("AR26").Value = 0
$ws.Range("AR29").Value = 1
$ws.Range("AR30").Value = 0
$ws.Range("AR31").Value = 0
$ws.Range("AR35").Value = 6
$ws.Range("AR36").Value = 0
$ws.Range("AR37").Value = 0
$ws.Range("AR38").Value = 0
$ws.Range("AR41").Value = 0
$ws.Range("AR42").Value = 0
$ws.Range("AR43").Value = 0

# Row 44 was also missing its week-40 (AQ) value; backfill it alongside
# the new week-41 (AR) value.
$ws.Range("AQ44").Value = 0
$ws.Range("AR44").Value = 0

$ws.Range("AR45").Value = 0
$ws.Range("AR46").Value = 0
$ws.Range("AR47").Value = 0
$ws.Range("AR48").Value = 0
$ws.Range("AR49").Value = 0
$ws.Range("AR50").Value = 0
$ws.Range("AR51").Value = 0
$ws.Range("AR53").Value = 0
$ws.Range("AR54").Value = 0
$ws.Range("AR55").Value = 0
$ws.Range("AR56").Value = 0
$ws.Range("AR57").Value = 0
$ws.Range("AR58").Value = 0
